{"js": "const replacements = [\n  [\"54\u00d734=1836\", \"21\u00d736=756\"],\n  [\"97\u00d720=1940\", \"27\u00d750=1350\"],\n  [\"51\u00d748=2448\", \"90\u00d729=2610\"],\n  [\"68\u00d749=3332\", \"30\u00d757=1710\"],\n  [\"57\u00d775=4275\", \"97\u00d749=4753\"],\n  [\"72\u00d753=3816\", \"69\u00d729=2001\"],\n  [\"69\u00d776=5244\", \"54\u00d746=2484\"],\n  [\"32\u00d766=2112\", \"40\u00d723=920\"],\n  [\"67\u00d731=2077\", \"32\u00d737=1184\"],\n  [\"48\u00d742=2016\", \"89\u00d715=1335\"],\n  [\"16\u00d749=784\", \"55\u00d759=3245\"],\n  [\"21\u00d785=1785\", \"77\u00d788=6776\"],\n  [\"84\u00d744=3696\", \"82\u00d730=2460\"],\n  [\"92\u00d754=4968\", \"69\u00d714=966\"],\n  [\"56\u00d718=1008\", \"30\u00d738=1140\"],\n  [\"60\u00d755=3300\", \"60\u00d791=5460\"],\n  [\"16\u00d754=864\", \"33\u00d728=924\"],\n  [\"31\u00d784=2604\", \"16\u00d758=928\"],\n  [\"79\u00d795=7505\", \"27\u00d726=702\"],\n  [\"14\u00d751=714\", \"32\u00d768=2176\"],\n  [\"19\u00d790=1710\", \"92\u00d755=5060\"],\n  [\"31\u00d771=2201\", \"72\u00d750=3600\"],\n  [\"53\u00d735=1855\", \"77\u00d714=1078\"],\n  [\"39\u00d762=2418\", \"67\u00d764=4288\"],\n  [\"35\u00d713=455\", \"57\u00d772=4104\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"54\u00d734=1836\"; New = \"21\u00d736=756\" },\n    @{ Old = \"97\u00d720=1940\"; New = \"27\u00d750=1350\" },\n    @{ Old = \"51\u00d748=2448\"; New = \"90\u00d729=2610\" },\n    @{ Old = \"68\u00d749=3332\"; New = \"30\u00d757=1710\" },\n    @{ Old = \"57\u00d775=4275\"; New = \"97\u00d749=4753\" },\n    @{ Old = \"72\u00d753=3816\"; New = \"69\u00d729=2001\" },\n    @{ Old = \"69\u00d776=5244\"; New = \"54\u00d746=2484\" },\n    @{ Old = \"32\u00d766=2112\"; New = \"40\u00d723=920\" },\n    @{ Old = \"67\u00d731=2077\"; New = \"32\u00d737=1184\" },\n    @{ Old = \"48\u00d742=2016\"; New = \"89\u00d715=1335\" },\n    @{ Old = \"16\u00d749=784\";  New = \"55\u00d759=3245\" },\n    @{ Old = \"21\u00d785=1785\"; New = \"77\u00d788=6776\" },\n    @{ Old = \"84\u00d744=3696\"; New = \"82\u00d730=2460\" },\n    @{ Old = \"92\u00d754=4968\"; New = \"69\u00d714=966\" },\n    @{ Old = \"56\u00d718=1008\"; New = \"30\u00d738=1140\" },\n    @{ Old = \"60\u00d755=3300\"; New = \"60\u00d791=5460\" },\n    @{ Old = \"16\u00d754=864\";  New = \"33\u00d728=924\" },\n    @{ Old = \"31\u00d784=2604\"; New = \"16\u00d758=928\" },\n    @{ Old = \"79\u00d795=7505\"; New = \"27\u00d726=702\" },\n    @{ Old = \"14\u00d751=714\";  New = \"32\u00d768=2176\" },\n    @{ Old = \"19\u00d790=1710\"; New = \"92\u00d755=5060\" },\n    @{ Old = \"31\u00d771=2201\"; New = \"72\u00d750=3600\" },\n    @{ Old = \"53\u00d735=1855\"; New = \"77\u00d714=1078\" },\n    @{ Old = \"39\u00d762=2418\"; New = \"67\u00d764=4288\" },\n    @{ Old = \"35\u00d713=455\";  New = \"57\u00d772=4104\" }\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $r.New, 2) | Out-Null\n}\n"}
